# Commit message: redefined prefix "ome" instead of ":" (base prefix name).
#
# The workbook lists RDF prefixes in the "@prefix" sheet (A=prefix name,
# B=namespace URI). Row 1 used to be the "base" prefix, stored with an
# empty prefix name (""), which made every "qname" in the other sheets
# that used the base prefix render as ":Something". This change gives the
# base prefix an explicit name, "ome", so those values become
# "ome:Something".

$wb = $excel.ActiveWorkbook

# --- 1. "@prefix" sheet: name the base prefix "ome" -------------------
$prefixSheet = $wb.Worksheets.Item("@prefix")
$prefixSheet.Range("A1").Value = "ome"

# --- 2. "Image" sheet: qnames using the base prefix --------------------
$imageSheet = $wb.Worksheets.Item("Image")
$imageSheet.Range("E3").Value = "ome:pixels"
$imageSheet.Range("B4").Value = "ome:Image"
$imageSheet.Range("E4").Value = "ome:Pixels"

# --- 3. "Pixels" sheet: qnames using the base prefix --------------------
$pixelsSheet = $wb.Worksheets.Item("Pixels")
$pixelsSheet.Range("D3").Value = "ome:pixelType"
$pixelsSheet.Range("E3").Value = "ome:dimensionOrder"
$pixelsSheet.Range("F3").Value = "ome:sizeC"
$pixelsSheet.Range("G3").Value = "ome:sizeT"
$pixelsSheet.Range("H3").Value = "ome:sizeX"
$pixelsSheet.Range("I3").Value = "ome:sizeY"
$pixelsSheet.Range("J3").Value = "ome:sizeZ"
$pixelsSheet.Range("K3").Value = "ome:binData"
$pixelsSheet.Range("B4").Value = "ome:Pixels"
$pixelsSheet.Range("D4").Value = "ome:PixelType"
$pixelsSheet.Range("E4").Value = "ome:DimensionOrder"
$pixelsSheet.Range("K4").Value = "ome:BinData"

# --- 4. "Binary_Data" sheet: qnames using the base prefix ---------------
$binSheet = $wb.Worksheets.Item("Binary_Data")
$binSheet.Range("C3").Value = "ome:bigEndian"
$binSheet.Range("D3").Value = "ome:compression"
$binSheet.Range("E3").Value = "ome:data"
$binSheet.Range("F3").Value = "ome:length"
$binSheet.Range("B4").Value = "ome:BinData"
$binSheet.Range("D4").Value = "ome:Compression"
